# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) on the per-job "profits" sheets after a pricing
# pull. Only the H:N metric columns move; leve identity columns (A:G) are
# untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3649.7856
$ws.Range("I62").Value = 3404.125
$ws.Range("K62").Value = 3404.125
$ws.Range("M62").Value = -2780.125

$ws.Range("H65").Value = 3649.7856
$ws.Range("I65").Value = 3404.125
$ws.Range("K65").Value = 17020.625
$ws.Range("M65").Value = -13900.625

$ws.Range("H76").Value = 8499.75
$ws.Range("J76").Value = 8499.75
$ws.Range("L76").Value = 8499.75
$ws.Range("N76").Value = -9129.75

$ws.Range("H79").Value = 8499.75
$ws.Range("J79").Value = 8499.75
$ws.Range("L79").Value = 8499.75
$ws.Range("N79").Value = -10683.75

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H92").Value = 803.3684
$ws.Range("I92").Value = 812.375
$ws.Range("J92").Value = 796.8182
$ws.Range("K92").Value = 812.375
$ws.Range("L92").Value = 796.8182
$ws.Range("M92").Value = 435.625
$ws.Range("N92").Value = -3292.8182

$ws.Range("H96").Value = 512.1111
$ws.Range("I96").Value = 569.875
$ws.Range("K96").Value = 1709.625
$ws.Range("M96").Value = -336.625

$ws.Range("H106").Value = 1969
$ws.Range("I106").Value = 1969
$ws.Range("K106").Value = 1969
$ws.Range("M106").Value = -1338

$ws.Range("H138").Value = 2804.7285
$ws.Range("I138").Value = 1481.3077
$ws.Range("J138").Value = 3586.75
$ws.Range("K138").Value = 4443.9231
$ws.Range("L138").Value = 10760.25
$ws.Range("M138").Value = 696.0769
$ws.Range("N138").Value = -21040.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22226240
$ws.Range("I32").Value = 23813580
$ws.Range("K32").Value = 23813580
$ws.Range("M32").Value = -23813293

$ws.Range("H44").Value = 39333
$ws.Range("J44").Value = 39333
$ws.Range("L44").Value = 39333
$ws.Range("N44").Value = -40309

$ws.Range("H51").Value = 50000
$ws.Range("I51").Value = 40000
$ws.Range("J51").Value = 52500
$ws.Range("K51").Value = 40000
$ws.Range("L51").Value = 52500
$ws.Range("M51").Value = -39244
$ws.Range("N51").Value = -54012

$ws.Range("H61").Value = 3131.25
$ws.Range("I61").Value = 3131.25
$ws.Range("K61").Value = 3131.25
$ws.Range("M61").Value = -2919.25

$ws.Range("H74").Value = 1320.6274
$ws.Range("I74").Value = 1380.2683
$ws.Range("J74").Value = 1076.1
$ws.Range("K74").Value = 1380.2683
$ws.Range("L74").Value = 1076.1
$ws.Range("M74").Value = -506.2683
$ws.Range("N74").Value = -2824.1

$ws.Range("H77").Value = 1320.6274
$ws.Range("I77").Value = 1380.2683
$ws.Range("J77").Value = 1076.1
$ws.Range("K77").Value = 6901.3415
$ws.Range("L77").Value = 5380.5
$ws.Range("M77").Value = -2533.3415
$ws.Range("N77").Value = -14116.5

$ws.Range("H88").Value = 2105.8333
$ws.Range("I88").Value = 2102
$ws.Range("J88").Value = 2107.75
$ws.Range("K88").Value = 2102
$ws.Range("L88").Value = 2107.75
$ws.Range("M88").Value = -1696
$ws.Range("N88").Value = -2919.75

$ws.Range("H91").Value = 2105.8333
$ws.Range("I91").Value = 2102
$ws.Range("J91").Value = 2107.75
$ws.Range("K91").Value = 2102
$ws.Range("L91").Value = 2107.75
$ws.Range("M91").Value = -698
$ws.Range("N91").Value = -4915.75

$ws.Range("H119").Value = 84799.60000000001
$ws.Range("J119").Value = 84799.60000000001
$ws.Range("L119").Value = 84799.60000000001
$ws.Range("N119").Value = -94475.60000000001

$ws.Range("H136").Value = 3131.25
$ws.Range("I136").Value = 3131.25
$ws.Range("K136").Value = 9393.75
$ws.Range("M136").Value = -6843.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 300
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -160
$ws.Range("N11").ClearContents()

$ws.Range("H20").Value = 2026.9546
$ws.Range("J20").Value = 3539.1428
$ws.Range("L20").Value = 3539.1428
$ws.Range("N20").Value = -4033.1428

$ws.Range("H86").Value = 4240
$ws.Range("I86").Value = 2716.5
$ws.Range("K86").Value = 2716.5
$ws.Range("M86").Value = -1593.5

$ws.Range("H89").Value = 4240
$ws.Range("I89").Value = 2716.5
$ws.Range("K89").Value = 13582.5
$ws.Range("M89").Value = -7966.5

$ws.Range("H107").Value = 7086.5356
$ws.Range("I107").Value = 4517.864
$ws.Range("K107").Value = 4517.864
$ws.Range("M107").Value = -2597.864

$ws.Range("H134").Value = 2394.7407
$ws.Range("I134").Value = 1823.25
$ws.Range("K134").Value = 5469.75
$ws.Range("M134").Value = -2934.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1993.7273
$ws.Range("I31").Value = 1850.5714
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1850.5714
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1555.5714
$ws.Range("N31").Value = -5590

$ws.Range("H34").Value = 1993.7273
$ws.Range("I34").Value = 1850.5714
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1850.5714
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1648.5714
$ws.Range("N34").Value = -5404

$ws.Range("H132").Value = 2810.0667
$ws.Range("I132").Value = 2225.3333
$ws.Range("K132").Value = 6675.999899999999
$ws.Range("M132").Value = -4145.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 22363.4
$ws.Range("J38").Value = 22363.4
$ws.Range("L38").Value = 22363.4
$ws.Range("N38").Value = -23289.4

$ws.Range("H52").Value = 38111.6
$ws.Range("I52").Value = 37030
$ws.Range("J52").Value = 38382
$ws.Range("K52").Value = 37030
$ws.Range("L52").Value = 38382
$ws.Range("M52").Value = -36771
$ws.Range("N52").Value = -38900

$ws.Range("H70").Value = 372537
$ws.Range("I70").Value = 557055.5
$ws.Range("K70").Value = 557055.5
$ws.Range("M70").Value = -556785.5

$ws.Range("H73").Value = 372537
$ws.Range("I73").Value = 557055.5
$ws.Range("K73").Value = 557055.5
$ws.Range("M73").Value = -556119.5

$ws.Range("H97").Value = 11618.444
$ws.Range("I97").Value = 509.42856
$ws.Range("K97").Value = 509.42856
$ws.Range("M97").Value = -13.42856

$ws.Range("H122").Value = 4911.3
$ws.Range("I122").Value = 4084.7693
$ws.Range("K122").Value = 12254.3079
$ws.Range("M122").Value = -9804.3079

$ws.Range("H132").Value = 4970.6313
$ws.Range("I132").Value = 5143.706
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 15431.118
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -12901.118
$ws.Range("N132").Value = -15558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 697
$ws.Range("J16").Value = 749.5
$ws.Range("L16").Value = 749.5
$ws.Range("N16").Value = -1089.5

$ws.Range("H132").Value = 2161.3674
$ws.Range("I132").Value = 1767.579
$ws.Range("K132").Value = 5302.737
$ws.Range("M132").Value = -2772.737

$ws.Range("H137").Value = 68214.28999999999
$ws.Range("J137").Value = 68846.16
$ws.Range("L137").Value = 68846.16
$ws.Range("N137").Value = -79046.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5449.75
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5449.75
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10899.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -13021.5

$ws.Range("H84").Value = 5449.75
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5449.75
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 54497.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -65105.5

$ws.Range("H107").Value = 832.625
$ws.Range("J107").Value = 865.6667
$ws.Range("L107").Value = 2597.0001
$ws.Range("N107").Value = -6437.0001

$ws.Range("H113").Value = 546.86957
$ws.Range("I113").Value = 512.44446
$ws.Range("J113").Value = 670.8
$ws.Range("K113").Value = 1537.33338
$ws.Range("L113").Value = 2012.4
$ws.Range("M113").Value = 632.66662
$ws.Range("N113").Value = -6352.4

$ws.Range("H122").Value = 1474.1904
$ws.Range("I122").Value = 1513.1177
$ws.Range("J122").Value = 1308.75
$ws.Range("K122").Value = 4539.3531
$ws.Range("L122").Value = 3926.25
$ws.Range("M122").Value = -2089.3531
$ws.Range("N122").Value = -8826.25
